# Auto-generated edit script: updates cryptocurrency price/volume data
# for rows 2-51 (Bitcoin ... PaxDollar) to match the scraped snapshot,
# including two pairs of rows that got re-ordered (13/14, 39/40, 41/42/43, 47/48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.711.95'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '1.764.35'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').Value = '''224.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').Value = '''0.544'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('D7').Value = '''1.01'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').Value = '''31.67'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('D9').Value = '''0.286'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').Value = '''0.0685'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').Value = '''0.0936'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('D12').Value = '2.029.16'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''10.88'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.35%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.760.54'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '33.816.31'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = '''0.611'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').Value = '''4.10'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = '''66.75'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '''237.49'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').Value = '0.0₃0769'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '''10.55'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '''4.02'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').Value = '''2.04'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').Value = '''158.94'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').Value = '''16.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').Value = '''6.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = '''0.0505'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').Value = '''3.59'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('D33').Value = '''3.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('D34').Value = '''1.77'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').Value = '1.374.27'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').Value = '''0.646'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').Value = '''1.03'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').Value = '''2.36'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''2.20'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.99%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '''2.65'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''0.896'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.30%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''76.86'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('D44').Value = '''13.26'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +12.94%  '
$ws.Range('D45').Value = '0.0₆0139'
$ws.Range('E45').Value = '  +16.28%  '
$ws.Range('E46').Value = '  +4.02%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').Value = '''0.0496'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''107.31'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('D49').Value = '''5.80'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').Value = '1.927.62'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('E51').Value = '  +0.64%  '
